$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1897668
$ws.Range("I38").Value = 1897668
$ws.Range("K38").Value = 5693004
$ws.Range("M38").Value = -5692632
$ws.Range("H39").Value = 1133890.4
$ws.Range("I39").Value = 1221058.9
$ws.Range("K39").Value = 3663176.7
$ws.Range("M39").Value = -3662880.7
$ws.Range("H64").Value = 34166.844
$ws.Range("I64").Value = 86032.664
$ws.Range("J64").Value = 3047.35
$ws.Range("K64").Value = 86032.664
$ws.Range("L64").Value = 3047.35
$ws.Range("M64").Value = -85784.664
$ws.Range("N64").Value = -3543.35
$ws.Range("H67").Value = 34166.844
$ws.Range("I67").Value = 86032.664
$ws.Range("J67").Value = 3047.35
$ws.Range("K67").Value = 86032.664
$ws.Range("L67").Value = 3047.35
$ws.Range("M67").Value = -85174.664
$ws.Range("N67").Value = -4763.35
$ws.Range("H76").Value = 4520.28
$ws.Range("J76").Value = 4880
$ws.Range("L76").Value = 4880
$ws.Range("N76").Value = -5510
$ws.Range("H79").Value = 4520.28
$ws.Range("J79").Value = 4880
$ws.Range("L79").Value = 4880
$ws.Range("N79").Value = -7064
$ws.Range("H98").Value = 931.6667
$ws.Range("I98").Value = 868
$ws.Range("J98").Value = 995.3333
$ws.Range("K98").Value = 868
$ws.Range("L98").Value = 995.3333
$ws.Range("M98").Value = 630
$ws.Range("N98").Value = -3991.3333
$ws.Range("H107").Value = 816.2
$ws.Range("I107").Value = 877.2353000000001
$ws.Range("J107").Value = 470.33334
$ws.Range("K107").Value = 877.2353000000001
$ws.Range("L107").Value = 470.33334
$ws.Range("M107").Value = 1042.7647
$ws.Range("N107").Value = -4310.33334
$ws.Range("H122").Value = 931.6667
$ws.Range("I122").Value = 868
$ws.Range("J122").Value = 995.3333
$ws.Range("K122").Value = 2604
$ws.Range("L122").Value = 2985.9999
$ws.Range("M122").Value = -154
$ws.Range("N122").Value = -7885.9999
$ws.Range("H132").Value = 5107181.5
$ws.Range("I132").Value = 5957707
$ws.Range("J132").Value = 4028.1428
$ws.Range("K132").Value = 17873121
$ws.Range("L132").Value = 12084.4284
$ws.Range("M132").Value = -17870591
$ws.Range("N132").Value = -17144.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23994.451
$ws.Range("I32").Value = 6227.3203
$ws.Range("J32").Value = 370453.5
$ws.Range("K32").Value = 6227.3203
$ws.Range("L32").Value = 370453.5
$ws.Range("M32").Value = -5940.3203
$ws.Range("N32").Value = -371027.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 22308.6
$ws.Range("J28").Value = 22308.6
$ws.Range("L28").Value = 22308.6
$ws.Range("N28").Value = -22798.6
$ws.Range("H31").Value = 29534.865
$ws.Range("I31").Value = 816.2727
$ws.Range("J31").Value = 50595.168
$ws.Range("K31").Value = 816.2727
$ws.Range("L31").Value = 50595.168
$ws.Range("M31").Value = -521.2727
$ws.Range("N31").Value = -51185.168
$ws.Range("H34").Value = 29534.865
$ws.Range("I34").Value = 816.2727
$ws.Range("J34").Value = 50595.168
$ws.Range("K34").Value = 816.2727
$ws.Range("L34").Value = 50595.168
$ws.Range("M34").Value = -614.2727
$ws.Range("N34").Value = -50999.168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 603.1429000000001
$ws.Range("I17").Value = 800
$ws.Range("J17").Value = 570.3333
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 1710.9999
$ws.Range("M17").Value = -2231
$ws.Range("N17").Value = -2048.9999
$ws.Range("H70").Value = 85834.5
$ws.Range("I70").Value = 168669
$ws.Range("K70").Value = 506007
$ws.Range("M70").Value = -505692
$ws.Range("H73").Value = 85834.5
$ws.Range("I73").Value = 168669
$ws.Range("K73").Value = 506007
$ws.Range("M73").Value = -504915
$ws.Range("H75").Value = 2579.2703
$ws.Range("I75").Value = 862.6
$ws.Range("J75").Value = 2847.5
$ws.Range("K75").Value = 2587.8
$ws.Range("L75").Value = 8542.5
$ws.Range("M75").Value = -1589.8
$ws.Range("N75").Value = -10538.5
$ws.Range("H78").Value = 2579.2703
$ws.Range("I78").Value = 862.6
$ws.Range("J78").Value = 2847.5
$ws.Range("K78").Value = 7763.400000000001
$ws.Range("L78").Value = 25627.5
$ws.Range("M78").Value = -2771.400000000001
$ws.Range("N78").Value = -35611.5
$ws.Range("H114").Value = 1127.9231
$ws.Range("I114").Value = 466
$ws.Range("J114").Value = 1541.625
$ws.Range("K114").Value = 1398
$ws.Range("L114").Value = 4624.875
$ws.Range("M114").Value = 1856
$ws.Range("N114").Value = -11132.875
$ws.Range("H131").Value = 842.48956
$ws.Range("J131").Value = 857.93475
$ws.Range("L131").Value = 2573.80425
$ws.Range("N131").Value = -12653.80425

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 111227220
$ws.Range("I80").Value = 143005860
$ws.Range("J80").Value = 2000
$ws.Range("K80").Value = 143005860
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = -143004862
$ws.Range("N80").Value = -3996
$ws.Range("H83").Value = 111227220
$ws.Range("I83").Value = 143005860
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 715029300
$ws.Range("L83").Value = 10000
$ws.Range("M83").Value = -715024308
$ws.Range("N83").Value = -19984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 4818.1665
$ws.Range("I68").Value = 2224.5
$ws.Range("J68").Value = 6115
$ws.Range("K68").Value = 2224.5
$ws.Range("L68").Value = 6115
$ws.Range("M68").Value = -1475.5
$ws.Range("N68").Value = -7613
$ws.Range("H71").Value = 4818.1665
$ws.Range("I71").Value = 2224.5
$ws.Range("J71").Value = 6115
$ws.Range("K71").Value = 11122.5
$ws.Range("L71").Value = 30575
$ws.Range("M71").Value = -7378.5
$ws.Range("N71").Value = -38063
$ws.Range("H82").Value = 1840.0625
$ws.Range("I82").Value = 1266
$ws.Range("J82").Value = 1972.5385
$ws.Range("K82").Value = 1266
$ws.Range("L82").Value = 1972.5385
$ws.Range("M82").Value = -905
$ws.Range("N82").Value = -2694.5385
$ws.Range("H85").Value = 1840.0625
$ws.Range("I85").Value = 1266
$ws.Range("J85").Value = 1972.5385
$ws.Range("K85").Value = 1266
$ws.Range("L85").Value = 1972.5385
$ws.Range("M85").Value = -18
$ws.Range("N85").Value = -4468.538500000001

Write-Output "Applied 165 cell updates"